$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 320
$ws.Range("I18").Value = 320
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 320
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("N18").Value = -36
$ws.Range("H32").Value = 919.5
$ws.Range("I32").Value = 489
$ws.Range("K32").Value = 489
$ws.Range("M32").Value = -163
$ws.Range("H100").Value = 1899.9131
$ws.Range("I100").Value = 1568.3158
$ws.Range("J100").Value = 3475
$ws.Range("K100").Value = 1568.3158
$ws.Range("L100").Value = 3475
$ws.Range("M100").Value = -1027.3158
$ws.Range("N100").Value = -4557
$ws.Range("H137").Value = 2504.25
$ws.Range("J137").Value = 4449.2
$ws.Range("L137").Value = 13347.6
$ws.Range("N137").Value = -18447.6
$ws.Range("H141").Value = 2943.818
$ws.Range("I141").Value = 2943.818
$ws.Range("K141").Value = 8831.454000000002
$ws.Range("M141").Value = -3651.454000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 40002384
$ws.Range("I61").Value = 76924410
$ws.Range("K61").Value = 76924410
$ws.Range("M61").Value = -76924198
$ws.Range("H74").Value = 58827484
$ws.Range("I74").Value = 71431304
$ws.Range("K74").Value = 71431304
$ws.Range("M74").Value = -71430430
$ws.Range("H77").Value = 58827484
$ws.Range("I77").Value = 71431304
$ws.Range("K77").Value = 357156520
$ws.Range("M77").Value = -357152152
$ws.Range("H88").Value = 102604.1
$ws.Range("I88").Value = 127380.125
$ws.Range("J88").Value = 3500
$ws.Range("K88").Value = 127380.125
$ws.Range("L88").Value = 3500
$ws.Range("M88").Value = -126974.125
$ws.Range("N88").Value = -4312
$ws.Range("H91").Value = 102604.1
$ws.Range("I91").Value = 127380.125
$ws.Range("J91").Value = 3500
$ws.Range("K91").Value = 127380.125
$ws.Range("L91").Value = 3500
$ws.Range("M91").Value = -125976.125
$ws.Range("N91").Value = -6308
$ws.Range("H97").Value = 843.0714
$ws.Range("I97").Value = 930.6111
$ws.Range("K97").Value = 930.6111
$ws.Range("M97").Value = -434.6111
$ws.Range("H102").Value = 6668722.5
$ws.Range("I102").Value = 14287247
$ws.Range("J102").Value = 2513.25
$ws.Range("K102").Value = 14287247
$ws.Range("L102").Value = 2513.25
$ws.Range("M102").Value = -14285625
$ws.Range("N102").Value = -5757.25
$ws.Range("H136").Value = 40002384
$ws.Range("I136").Value = 76924410
$ws.Range("K136").Value = 230773230
$ws.Range("M136").Value = -230770680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2175.8333
$ws.Range("I86").Value = 2329.375
$ws.Range("J86").Value = 947.5
$ws.Range("K86").Value = 2329.375
$ws.Range("L86").Value = 947.5
$ws.Range("M86").Value = -1206.375
$ws.Range("N86").Value = -3193.5
$ws.Range("H89").Value = 2175.8333
$ws.Range("I89").Value = 2329.375
$ws.Range("J89").Value = 947.5
$ws.Range("K89").Value = 11646.875
$ws.Range("L89").Value = 4737.5
$ws.Range("M89").Value = -6030.875
$ws.Range("N89").Value = -15969.5
$ws.Range("H134").Value = 23814090
$ws.Range("I134").Value = 25004594
$ws.Range("K134").Value = 75013782
$ws.Range("M134").Value = -75011247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 5400
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = $null
$ws.Range("H86").Value = 12738.846
$ws.Range("I86").Value = 9351
$ws.Range("K86").Value = 9351
$ws.Range("M86").Value = -8228
$ws.Range("H89").Value = 12738.846
$ws.Range("I89").Value = 9351
$ws.Range("K89").Value = 46755
$ws.Range("M89").Value = -41139
$ws.Range("H132").Value = 200001630
$ws.Range("I132").Value = 200001630
$ws.Range("K132").Value = 600004890
$ws.Range("M132").Value = -600002360
$ws.Range("H133").Value = 96246
$ws.Range("J133").Value = 96246
$ws.Range("L133").Value = 96246
$ws.Range("N133").Value = -101306

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 53.875
$ws.Range("I2").Value = 50
$ws.Range("J2").Value = 55.166668
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 331.000008
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -557.000008
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 47.333332
$ws.Range("J23").Value = 52
$ws.Range("K23").Value = 141.999996
$ws.Range("L23").Value = 156
$ws.Range("M23").Value = 93.00000399999999
$ws.Range("N23").Value = -626
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = $null
$ws.Range("N48").Value = 0
$ws.Range("H55").Value = 1316.6666
$ws.Range("H86").Value = 755.26666
$ws.Range("I86").Value = 537.3333
$ws.Range("K86").Value = 1611.9999
$ws.Range("M86").Value = -425.9999
$ws.Range("H89").Value = 755.26666
$ws.Range("I89").Value = 537.3333
$ws.Range("K89").Value = 4835.9997
$ws.Range("M89").Value = 1092.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 24932
$ws.Range("J104").Value = 24932
$ws.Range("L104").Value = 24932
$ws.Range("N104").Value = -31920
$ws.Range("H113").Value = 56084.65
$ws.Range("I113").Value = 65305.53
$ws.Range("J113").Value = 3833
$ws.Range("K113").Value = 65305.53
$ws.Range("L113").Value = 3833
$ws.Range("M113").Value = -63135.53
$ws.Range("N113").Value = -8173
$ws.Range("H132").Value = 12502709
$ws.Range("I132").Value = 13891577
$ws.Range("K132").Value = 41674731
$ws.Range("M132").Value = -41672201

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3948.125
$ws.Range("I40").Value = 3798
$ws.Range("K40").Value = 3798
$ws.Range("M40").Value = -3662
$ws.Range("H55").Value = 453.875
$ws.Range("I55").Value = 332.875
$ws.Range("J55").Value = 574.875
$ws.Range("K55").Value = 332.875
$ws.Range("L55").Value = 574.875
$ws.Range("M55").Value = -159.875
$ws.Range("N55").Value = -920.875
$ws.Range("H82").Value = 1727.8889
$ws.Range("I82").Value = 1742.8948
$ws.Range("K82").Value = 1742.8948
$ws.Range("M82").Value = -1381.8948
$ws.Range("H85").Value = 1727.8889
$ws.Range("I85").Value = 1742.8948
$ws.Range("K85").Value = 1742.8948
$ws.Range("M85").Value = -494.8948
$ws.Range("H93").Value = 2767.5557
$ws.Range("I93").Value = 1648
$ws.Range("K93").Value = 1648
$ws.Range("M93").Value = -400
$ws.Range("H94").Value = 15000
$ws.Range("J94").Value = 15000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
$ws.Range("H100").Value = 8681190
$ws.Range("I100").Value = 9982798
$ws.Range("J100").Value = 3799.6667
$ws.Range("K100").Value = 9982798
$ws.Range("L100").Value = 3799.6667
$ws.Range("M100").Value = -9982257
$ws.Range("N100").Value = -4881.6667
$ws.Range("H132").Value = 43640470
$ws.Range("I132").Value = 48003820
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 144011460
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -144008930
$ws.Range("N132").Value = -26060
$ws.Range("H136").Value = 1986
$ws.Range("I136").Value = 1986
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5958
$ws.Range("L136").Value = $null
$ws.Range("M136").Value = -3408
$ws.Range("N136").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 44906.13
$ws.Range("I81").Value = 51113.05
$ws.Range("K81").Value = 102226.1
$ws.Range("M81").Value = -101165.1
$ws.Range("H84").Value = 44906.13
$ws.Range("I84").Value = 51113.05
$ws.Range("K84").Value = 511130.5
$ws.Range("M84").Value = -505826.5
$ws.Range("H100").Value = 2048.7273
$ws.Range("I100").Value = 1953.6
$ws.Range("K100").Value = 3907.2
$ws.Range("M100").Value = -3366.2
$ws.Range("H122").Value = 1712.5454
$ws.Range("I122").Value = 1683.8
$ws.Range("K122").Value = 5051.4
$ws.Range("M122").Value = -2601.4
$ws.Range("H132").Value = 12508206
$ws.Range("I132").Value = 17864022
$ws.Range("K132").Value = 53592066
$ws.Range("M132").Value = -53589536
$ws.Range("H136").Value = 26317808
$ws.Range("I136").Value = 31252072
$ws.Range("K136").Value = 93756216
$ws.Range("M136").Value = -93753666
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = $null
$ws.Range("N141").Value = 0
